$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 is the last existing data row; copy its cell structure/format down to the
# new row 5 so the new row inherits the same per-cell types (including the blank
# inline-string cell in column AG) and the bold/centered/bordered style used on
# column A.
$ws.Range("A4:DK4").Copy($ws.Range("A5:DK5"))

$values = @("2021年", 12, 36.5, -32.5, 44.1, 2.2, -21.2, 18.1, 7.9, 9, -6.6, 5.2, 19.2, -37.5, 0.2, 11.9, -13.7, 62.1, 32.4, 14.8, 17.8, 13.9, 26.9, 22.1, 24.4, 62.4, 19.9, 22.5, 27.6, 10.6, 16.7, -40.4, $null, -18.2, -6.3, -13.4, 8.6, 8.4, 7.3, -6, -6.6, -22, 15.9, 4.7, 68.7, 35.6, 50.2, 7.5, 21, -0.3, 0.2, 11.7, 5.6, 11.5, 18.1, 39.8, 16.3, 4, 23.5, 6, 7.2, 20.3, 22.1, 0.4, 3, -1, 2.3, 11.2, 14.4, 7.8, 1.5, -0.6, -5.4, 3.5, 6.2, 24.4, 26.1, 11.4, 28.5, 3.5, 14.3, -23, 1.8, 17.5, 22.9, 19.3, 16.7, 4.5, 26.6, 10.6, -16.8, 16.9, 17.9, 32.5, 2.2, -36.7, 27, 16.9, 12.2, -0.6, 10.4, 22.1, 10.1, 52.5, 19, -0.3, 28, 2.8, 0.6, 18.6, 24.2, 15.7, 17.1, 35.3, 33.2)

for ($i = 0; $i -lt $values.Length; $i++) {
    if ($null -ne $values[$i]) {
        $ws.Cells.Item(5, $i + 1).Value2 = $values[$i]
    }
}
